$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary counters (row 13) ---
# Cant. Trabajadores: 3 -> 2
$ws.Range("C13").Value = 2
# Cant. Periodos: 2 -> 3
$ws.Range("F13").Value = 3

# --- Update worker data table (rows 16-18) ---
# Row 16: EUCLIDIA HERAZO SIMARRA -> ADEMIR ENRIQUE PADILLA BORNACHERA
$ws.Range("C16").Value = "73228947"
$ws.Range("D16").Value = "ADEMIR ENRIQUE PADILLA BORNACHERA"
$ws.Range("E16").Value = "2402"
$ws.Range("F16").Value = 59560
$ws.Range("G16").Value = 1623900

# Row 17: KAREN MARGARITA LARA ARROYO stays, values unchanged
$ws.Range("C17").Value = "1050034108"
$ws.Range("D17").Value = "KAREN MARGARITA LARA ARROYO"
$ws.Range("E17").Value = "2507"
$ws.Range("F17").Value = 75884
$ws.Range("G17").Value = 1897100

# Row 18: new entry - KAREN MARGARITA LARA ARROYO / period 2508
$ws.Range("C18").Value = "1050034108"
$ws.Range("D18").Value = "KAREN MARGARITA LARA ARROYO"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 2529
$ws.Range("G18").Value = 1897100
